$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# New "TIME mode" configure template rows appended to the Translation sheet
# (rows 6-10), following the same layout as the existing rows 4-5:
#   B = Text Id, C = Typography Name, D = Alignment, E = Direction, F = Text

$rows = @(
    @("SingleUseId3", "Default", "Left", "LTR", "TIME INTERVAL/MODE: INDEPENDENT"),
    @("SingleUseId4", "Default", "Left", "LTR", "INPUT"),
    @("SingleUseId5", "Default", "Left", "LTR", "CLOSK"),
    @("SingleUseId6", "Default", "Left", "LTR", "TI SETUP"),
    @("SingleUseId7", "Default", "Left", "LTR", "SESSION SETUP")
)

$r = 6
foreach ($row in $rows) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $r = $r + 1
}
